$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68, shifting existing rows 68-154 down to 69-155.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new record's data.
$ws.Range("A68").Value = 1
$ws.Range("B68").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C68").Value = 'Arica y Parinacota'
$ws.Range("D68").Value = (Get-Date -Year 2023 -Month 4 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E68").Value = 15
$ws.Range("F68").Value = 'Fruta'
$ws.Range("G68").Value = 100102
$ws.Range("H68").Value = 'Cítricos'
$ws.Range("I68").Value = 100102004
$ws.Range("J68").Value = 'Mandarina'
$ws.Range("K68").Value = 'Murcott'
$ws.Range("L68").Value = 'Segunda'
$ws.Range("M68").Value = 250
$ws.Range("N68").Value = 22000
$ws.Range("O68").Value = 23000
$ws.Range("P68").Value = 22600
$ws.Range("Q68").Value = '$/caja 20 kilos'
$ws.Range("R68").Value = "Región de O'Higgins"
$ws.Range("S68").Value = 1130
$ws.Range("T68").Value = 20
